$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 201 (shifts existing rows 201:299 down to 202:300)
$ws.Rows("201:201").Insert()

# Populate the newly inserted row 201 with a new data record
$ws.Cells.Item(201, 1).Value = 7
$ws.Cells.Item(201, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(201, 3).Value = "Ñuble"
$ws.Cells.Item(201, 4).Value = 44900
$ws.Cells.Item(201, 5).Value = 16
$ws.Cells.Item(201, 6).Value = 100112006
$ws.Cells.Item(201, 7).Value = "Repollo"
$ws.Cells.Item(201, 8).Value = "Crespo record"
$ws.Cells.Item(201, 9).Value = "Primera"
$ws.Cells.Item(201, 10).Value = 500
$ws.Cells.Item(201, 11).Value = 1500
$ws.Cells.Item(201, 12).Value = 1600
$ws.Cells.Item(201, 13).Value = 1550
$ws.Cells.Item(201, 14).Value = "$/unidad"
$ws.Cells.Item(201, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(201, 16).Value = 1550
$ws.Cells.Item(201, 17).Value = 1
$ws.Cells.Item(201, 18).Value = "Hortaliza"

# Keep the date cell's number format consistent with the other date cells in column D
$ws.Cells.Item(201, 4).NumberFormat = $ws.Cells.Item(202, 4).NumberFormat
